# fall 24 week 4 inputs
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("C2").Value = 10.19

$ws.Range("B3").Value = 9.81
$ws.Range("D3").Value = 10.12

$ws.Range("C4").Value = 9.880000000000001
$ws.Range("E4").Value = 10.64

$ws.Range("D5").Value = 9.359999999999999
$ws.Range("F5").Value = 10.19
$ws.Range("G5").Value = 9.710000000000001

$ws.Range("E6").Value = 9.81
$ws.Range("G6").Value = 10.48
$ws.Range("H6").Value = 10.53

$ws.Range("E7").Value = 10.29
$ws.Range("F7").Value = 9.52
$ws.Range("H7").Value = 9.91
$ws.Range("J7").Value = 9.67

$ws.Range("F8").Value = 9.470000000000001
$ws.Range("G8").Value = 10.09

$ws.Range("G10").Value = 10.33
